$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Title"
$ws.Range("C1").Value = "Text"
$ws.Range("D1").Value = "Author"
$ws.Range("F1").Value = "Result"

$ws.Range("F6").Select()
